$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values (back-up data from API): E6, E8 changed and E9 added
$ws.Range("E6").Value = 0.9
$ws.Range("E8").Value = 0.9
$ws.Range("E4").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").Value = 0.8

# Update the active selection to reflect where the author last clicked
$ws.Range("E10").Select()
